# Apply Leve profit-table value updates (scheduled price refresh).
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 98 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 991.5417
$ws.Range("I98").Value = 561.55554
$ws.Range("J98").Value = 2281.5
$ws.Range("K98").Value = 561.55554
$ws.Range("L98").Value = 2281.5
$ws.Range("M98").Value = 936.44446
$ws.Range("N98").Value = -5277.5

# Sheet ALC, row 122 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 991.5417
$ws.Range("I122").Value = 561.55554
$ws.Range("J122").Value = 2281.5
$ws.Range("K122").Value = 1684.66662
$ws.Range("L122").Value = 6844.5
$ws.Range("M122").Value = 765.33338
$ws.Range("N122").Value = -11744.5

# Sheet ALC, row 134 (Leve Item ID 41997)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 55000
$ws.Range("J134").Value = 55000
$ws.Range("L134").Value = 55000
$ws.Range("N134").Value = -65140

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1362.4482
$ws.Range("I137").Value = 1167.0741
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 3501.2223
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -951.2223000000004
$ws.Range("N137").Value = -17100

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1340.3334
$ws.Range("I61").Value = 1296.7222
$ws.Range("J61").Value = 1602
$ws.Range("K61").Value = 1296.7222
$ws.Range("L61").Value = 1602
$ws.Range("M61").Value = -1084.7222
$ws.Range("N61").Value = -2026

# Sheet ARM, row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 911.325
$ws.Range("I74").Value = 974.3
$ws.Range("J74").Value = 722.4
$ws.Range("K74").Value = 974.3
$ws.Range("L74").Value = 722.4
$ws.Range("M74").Value = -100.3
$ws.Range("N74").Value = -2470.4

# Sheet ARM, row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 911.325
$ws.Range("I77").Value = 974.3
$ws.Range("J77").Value = 722.4
$ws.Range("K77").Value = 4871.5
$ws.Range("L77").Value = 3612
$ws.Range("M77").Value = -503.5
$ws.Range("N77").Value = -12348

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 164063.33
$ws.Range("I132").Value = 218620.66
$ws.Range("J132").Value = 7211
$ws.Range("K132").Value = 655861.98
$ws.Range("L132").Value = 21633
$ws.Range("M132").Value = -653331.98
$ws.Range("N132").Value = -26693

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1340.3334
$ws.Range("I136").Value = 1296.7222
$ws.Range("J136").Value = 1602
$ws.Range("K136").Value = 3890.1666
$ws.Range("L136").Value = 4806
$ws.Range("M136").Value = -1340.1666
$ws.Range("N136").Value = -9906

# Sheet BSM, row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 203528.14
$ws.Range("I134").Value = 217744.42
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 653233.26
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -650698.26
$ws.Range("N134").Value = -18570

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2083.394
$ws.Range("I31").Value = 1589.6364
$ws.Range("J31").Value = 3070.9092
$ws.Range("K31").Value = 1589.6364
$ws.Range("L31").Value = 3070.9092
$ws.Range("M31").Value = -1294.6364
$ws.Range("N31").Value = -3660.9092

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2083.394
$ws.Range("I34").Value = 1589.6364
$ws.Range("J34").Value = 3070.9092
$ws.Range("K34").Value = 1589.6364
$ws.Range("L34").Value = 3070.9092
$ws.Range("M34").Value = -1387.6364
$ws.Range("N34").Value = -3474.9092

# Sheet CRP, row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1235.4375
$ws.Range("I58").Value = 1135.6552
$ws.Range("J58").Value = 2200
$ws.Range("K58").Value = 1135.6552
$ws.Range("L58").Value = 2200
$ws.Range("M58").Value = -932.6551999999999
$ws.Range("N58").Value = -2606

# Sheet CRP, row 99 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1733.2222
$ws.Range("I99").Value = 1720
$ws.Range("J99").Value = 1749.75
$ws.Range("K99").Value = 1720
$ws.Range("L99").Value = 1749.75
$ws.Range("M99").Value = -222
$ws.Range("N99").Value = -4745.75

# Sheet CRP, row 107 (Leve Item ID 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1410.1875
$ws.Range("I107").Value = 1535
$ws.Range("J107").Value = 1202.1666
$ws.Range("K107").Value = 1535
$ws.Range("L107").Value = 1202.1666
$ws.Range("M107").Value = 385
$ws.Range("N107").Value = -5042.1666

# Sheet CRP, row 126 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1733.2222
$ws.Range("I126").Value = 1720
$ws.Range("J126").Value = 1749.75
$ws.Range("K126").Value = 5160
$ws.Range("L126").Value = 5249.25
$ws.Range("M126").Value = -2690
$ws.Range("N126").Value = -10189.25

# Sheet CRP, row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1998.0358
$ws.Range("I132").Value = 1595.3182
$ws.Range("J132").Value = 3474.6667
$ws.Range("K132").Value = 4785.9546
$ws.Range("L132").Value = 10424.0001
$ws.Range("M132").Value = -2255.9546
$ws.Range("N132").Value = -15484.0001

# Sheet CRP, row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4835.0605
$ws.Range("I134").Value = 5474.6294
$ws.Range("J134").Value = 1957
$ws.Range("K134").Value = 16423.8882
$ws.Range("L134").Value = 5871
$ws.Range("M134").Value = -13888.8882
$ws.Range("N134").Value = -10941

# Sheet CRP, row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1235.4375
$ws.Range("I136").Value = 1135.6552
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 3406.9656
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -856.9655999999995
$ws.Range("N136").Value = -11700

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2447.8823
$ws.Range("I122").Value = 2493.6667
$ws.Range("J122").Value = 2338
$ws.Range("K122").Value = 7481.000100000001
$ws.Range("L122").Value = 7014
$ws.Range("M122").Value = -5031.000100000001
$ws.Range("N122").Value = -11914

# Sheet GSM, row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1909.2059
$ws.Range("I132").Value = 1675.0333
$ws.Range("J132").Value = 3665.5
$ws.Range("K132").Value = 5025.0999
$ws.Range("L132").Value = 10996.5
$ws.Range("M132").Value = -2495.0999
$ws.Range("N132").Value = -16056.5

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2857.4167
$ws.Range("I132").Value = 2128.2354
$ws.Range("K132").Value = 6384.706200000001
$ws.Range("M132").Value = -3854.706200000001

# Sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1647.8518
$ws.Range("I136").Value = 1605.8948
$ws.Range("J136").Value = 1747.5
$ws.Range("K136").Value = 4817.6844
$ws.Range("L136").Value = 5242.5
$ws.Range("M136").Value = -2267.6844
$ws.Range("N136").Value = -10342.5

# Sheet WVR, row 125 (Leve Item ID 34276)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 53167.855
$ws.Range("J125").Value = 53167.855
$ws.Range("L125").Value = 53167.855
$ws.Range("N125").Value = -63007.855

# Sheet WVR, row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1257.1154
$ws.Range("I126").Value = 767.2727
$ws.Range("J126").Value = 3951.25
$ws.Range("K126").Value = 2301.8181
$ws.Range("L126").Value = 11853.75
$ws.Range("M126").Value = 168.1819
$ws.Range("N126").Value = -16793.75

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3879.8064
$ws.Range("I132").Value = 3758.0454
$ws.Range("J132").Value = 4177.4443
$ws.Range("K132").Value = 11274.1362
$ws.Range("L132").Value = 12532.3329
$ws.Range("M132").Value = -8744.136200000001
$ws.Range("N132").Value = -17592.3329

# Sheet WVR, row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1434.5946
$ws.Range("I136").Value = 1336.6666
$ws.Range("J136").Value = 1854.2858
$ws.Range("K136").Value = 4009.9998
$ws.Range("L136").Value = 5562.857400000001
$ws.Range("M136").Value = -1459.9998
$ws.Range("N136").Value = -10662.8574
